# Applies the "Automatic update of files" edit:
#   1. Updates the "Förändrad" (changed) date in column C from 45184 to 45186
#      for every data row.
#   2. Adds a friendly display name (the report's "Beteckning", column A) as
#      the second argument of every HYPERLINK() formula found in columns
#      S, T, V, W, X and Y, e.g.
#        =HYPERLINK("...A 412-2023.xlsx")
#      becomes
#        =HYPERLINK("...A 412-2023.xlsx", "A 412-2023")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = 1   # Beteckning
$colC = 3   # Förändrad (the date column that gets bumped)

# HYPERLINK columns that may need a friendly-name second argument.
$linkCols = @(19, 20, 22, 23, 24, 25)   # S, T, V, W, X, Y

$oldDate = 45184
$newDate = 45186

# Cell values that are formatted as dates come back as System.DateTime
# instead of a raw serial number; normalise to the Excel serial (OA date).
function Get-SerialValue($val) {
    if ($null -eq $val) { return $null }
    if ($val.GetType().FullName -eq "System.DateTime") {
        return $val.ToOADate()
    }
    return $val
}

# Find the last used data row by walking up from the bottom of column A.
$lastRow = $ws.Cells.Item(1048576, $colA).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {

    # --- 1) Bump the "Förändrad" date in column C ---
    $cCell = $ws.Cells.Item($r, $colC)
    $cVal = Get-SerialValue ($cCell.Value())
    if ($cVal -eq $oldDate) {
        $cCell.Value = $newDate
    }

    # --- 2) Add the friendly name to any HYPERLINK formulas in this row ---
    $label = $ws.Cells.Item($r, $colA).Value()
    # Escape any embedded double quotes the way Excel formulas expect ("" ).
    $escapedLabel = [string]$label -replace '"', '""'

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($r, $col)
        if ($cell.HasFormula) {
            $formula = $cell.Formula
            $body = $formula
            if ($body.StartsWith("=")) { $body = $body.Substring(1) }

            if ($body -match "^HYPERLINK\(" -and -not ($body.Contains(","))) {
                $lastParen = $body.LastIndexOf(")")
                $newBody = $body.Substring(0, $lastParen) + ', "' + $escapedLabel + '"' + $body.Substring($lastParen)
                $cell.Formula = "=" + $newBody
            }
        }
    }
}
